$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")
$ws.Activate()
try {
  $excel.ActiveWindow.TopLeftCell = $ws.Range("A19")
  Write-Host "set topleft via property"
} catch {
  Write-Host "ERR1: $_"
}
